$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-01-10 Wednesday"; new = "2024-01-11 Thursday"},
    @{old = "574×7=4018"; new = "307×4=1228"},
    @{old = "459×7=3213"; new = "402×7=2814"},
    @{old = "385×7=2695"; new = "344×2=688"},
    @{old = "134×8=1072"; new = "618×7=4326"},
    @{old = "445×6=2670"; new = "239×6=1434"},
    @{old = "198×6=1188"; new = "185×4=740"},
    @{old = "237×7=1659"; new = "732×8=5856"},
    @{old = "822×6=4932"; new = "647×8=5176"},
    @{old = "599×8=4792"; new = "196×4=784"},
    @{old = "806×6=4836"; new = "667×3=2001"},
    @{old = "740×5=3700"; new = "406×4=1624"},
    @{old = "649×6=3894"; new = "229×8=1832"},
    @{old = "719×6=4314"; new = "424×3=1272"},
    @{old = "303×5=1515"; new = "541×9=4869"},
    @{old = "296×8=2368"; new = "511×3=1533"},
    @{old = "319×7=2233"; new = "305×6=1830"},
    @{old = "696×5=3480"; new = "230×8=1840"},
    @{old = "728×2=1456"; new = "964×2=1928"},
    @{old = "548×8=4384"; new = "508×5=2540"},
    @{old = "343×8=2744"; new = "624×9=5616"},
    @{old = "604×3=1812"; new = "693×8=5544"},
    @{old = "788×7=5516"; new = "844×4=3376"},
    @{old = "894×9=8046"; new = "304×9=2736"},
    @{old = "356×5=1780"; new = "895×8=7160"},
    @{old = "527×7=3689"; new = "597×3=1791"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
